# Applies the commit:
#  "Added fom costs and adjusted mapping; also adjusted has_state to
#   achieve 'true' instead of 'True' to avoid read in errors"
#
# This:
#  1. Removes the now-redundant "fom_cost" column from the "Nodes" sheet
#     (columns D..H shift left to become C..G).
#  2. Re-maps/re-orders the node rows both on the "Definition" sheet
#     (column A only) and on the "Nodes" sheet (whole rows), to the new
#     canonical node order.
#  3. Re-writes the has_state values as the literal text "true" (instead
#     of the Excel boolean TRUE which used to serialize as "True").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Definition sheet: re-order the node names in column A (rows 7-18).
#    Category (B) and fom_cost (C) stay blank/"node" for all these rows,
#    so only the names need to move.
# ---------------------------------------------------------------------
$wsDef = $wb.Worksheets.Item("Definition")

$defNames = @(
    "Power_Kasso",
    "E-Methanol_Kasso",
    "E-Methanol_storage_Kasso",
    "Vaporized_Carbon_Dioxide",
    "Waste_Heat",
    "Carbon_Dioxide",
    "Hydrogen_Kasso",
    "Raw_Methanol",
    "District_Heating",
    "Water",
    "Hydrogen_storage_Kasso",
    "Power_Wholesale"
)

$defArr = New-Object 'object[,]' $defNames.Length,1
for ($i = 0; $i -lt $defNames.Length; $i++) {
    $defArr[$i,0] = $defNames[$i]
}
$wsDef.Range("A7:A18").Value = $defArr

# ---------------------------------------------------------------------
# 2) Nodes sheet: drop the "fom_cost" column (column C) so that
#    balance_type / has_state / node_state_cap / frac_state_loss /
#    node_slack_penalty shift from D-H to C-G.
# ---------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("Nodes")
$wsNodes.Columns.Item(3).Delete()

# ---------------------------------------------------------------------
# 3) Nodes sheet: re-order the rows to the new node order, carrying each
#    node's (now shifted) data along with it, and write has_state as the
#    literal text "true" (leading "'" forces text instead of boolean).
# ---------------------------------------------------------------------
$nodeRows = @(
    @("Power_Kasso",               "node", "balance_type_node", "",      "",     "",  100000),
    @("E-Methanol_Kasso",          "node", "balance_type_node", "",      "",     "",  100000),
    @("E-Methanol_storage_Kasso",  "node", "balance_type_node", "'true", 100000, 0,   100000),
    @("Vaporized_Carbon_Dioxide",  "node", "balance_type_node", "",      "",     "",  100000),
    @("Waste_Heat",                "node", "balance_type_node", "",      "",     "",  ""),
    @("Carbon_Dioxide",            "node", "balance_type_none", "",      "",     "",  ""),
    @("Hydrogen_Kasso",            "node", "balance_type_node", "",      "",     "",  100000),
    @("Raw_Methanol",              "node", "balance_type_node", "",      "",     "",  100000),
    @("District_Heating",          "node", "balance_type_none", "",      "",     "",  ""),
    @("Water",                     "node", "balance_type_none", "",      "",     "",  ""),
    @("Hydrogen_storage_Kasso",    "node", "balance_type_node", "'true", 100000, 0,   100000),
    @("Power_Wholesale",           "node", "balance_type_none", "",      "",     "",  "")
)

$nodeArr = New-Object 'object[,]' $nodeRows.Length,7
for ($r = 0; $r -lt $nodeRows.Length; $r++) {
    $row = $nodeRows[$r]
    for ($c = 0; $c -lt 7; $c++) {
        $nodeArr[$r,$c] = $row[$c]
    }
}
$wsNodes.Range("A2:G13").Value = $nodeArr
